# Auto-generated script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '58.995.26'
$ws.Range("E2").Value = '  -6.02%  '
$ws.Range("D3").Value = '2.443.13'
$ws.Range("E3").Value = '  -8.79%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -2.74%  '
$ws.Range("D9").Value = '2.457.45'
$ws.Range("E9").Value = '  -8.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0990'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.159'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("E13").Value = '  -4.37%  '
$ws.Range("D14").Value = '2.893.02'
$ws.Range("E14").Value = '  -8.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.60%  '
$ws.Range("D16").Value = '58.935.86'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000137'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.42%  '
$ws.Range("D18").Value = '2.504.05'
$ws.Range("E18").Value = '  -6.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.448'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -11.94%  '
$ws.Range("E26").Value = '  -4.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.977'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.70%  '
$ws.Range("E30").Value = '  -5.70%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0₃0771'
$ws.Range("E31").Value = '  -8.98%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.997'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '156.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.43'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.59%  '
$ws.Range("E38").Value = '  -4.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.80'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '313.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.831'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.25%  '
$ws.Range("E43").Value = '  -6.78%  '
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("E46").Value = '  -4.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0933'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0525'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0229'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.56%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.43%  '

$wb.Save()
